# Updating group diary 24/10
# Adds the 24 Oct 2023 meeting entry (row 10) to Sheet2 of the group diary
# and moves the sheet's active-cell selection to H5, matching the author's
# edit on their machine.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New meeting row: date, start/end time, members present, discussion notes.
$ws.Range("A10").Value = 45223                      # 24/10/2023
$ws.Range("B10").Value = 0.39583333333333331         # 09:30
$ws.Range("C10").Value = 0.41666666666666669         # 10:00
$ws.Range("D10").Value = "All"
$ws.Range("E10").Value = "Checking on the working process of groups member and discussing about the presentations and adjustments. Assigning editing tasks."

# The long discussion note wraps across multiple lines, so the row grows taller.
$ws.Rows.Item(10).RowHeight = 51

# Leave the cursor where the author left it when they saved.
[void]$ws.Range("H5").Select()
